$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column R, rows 4-14 (same pattern/styles as column Q in each row)
$values = @{
    4  = 2020
    5  = 2.1
    6  = 2.4
    7  = 1.4
    8  = 3.2
    9  = 2.4
    10 = 0.8
    11 = 2.2000000000000002
    12 = 4.5
    13 = 1.4
    14 = 3.2
}

foreach ($row in 4..14) {
    $qCell = $ws.Range("Q$row")
    $rCell = $ws.Range("R$row")

    # Copy formatting (style) from Q to R so the new column matches the table's look
    $qCell.Copy()
    $rCell.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false

    # Set the cell's value
    $rCell.Value = $values[$row]
}

# Update the selection to match the saved state in the diff
$ws.Range("R16:R17").Select()
